# Complejidad.xlsx — add the two new BinomioDeNewton method rows
# (obtenerCoeficiente / obtenerPolinomio) to the "Hoja1" complexity table,
# and move the sheet's selection down onto the newly filled rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 8: obtenerCoeficiente -> n -------------------------------------
$ws.Range("A8").Value = "obtenerCoeficiente"
$ws.Range("B8").Value = "n"

# --- Row 9: obtenerPolinomio -> n^2 -------------------------------------
$ws.Range("A9").Value = "obtenerPolinomio"
$ws.Range("B9").Value = "n^2"

# Both new complexity cells are centered, like the rest of column B.
$ws.Range("B8:B9").HorizontalAlignment = -4108   # xlCenter

# Move the active selection onto the rows that were just filled in.
$ws.Range("B8:B9").Select()
